$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New collaborator rows (4, 5, 6) ---

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Ricardo Avila"
$ws.Range("C5").Value = "Comercial"
$ws.Range("D5").Value = 573163215029
$ws.Range("D5").NumberFormat = "0"
$ws.Range("E5").Value = "ravila@suraelec.com"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Celina Charris"
$ws.Range("C6").Value = "Lider Regional"
$ws.Range("D6").Value = 573183117195
$ws.Range("D6").NumberFormat = "0"
$ws.Range("E6").Value = "ccharris@suraelec.com"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Juan J Lozano"
$ws.Range("C7").Value = "Gerente"
$ws.Range("D7").Value = 573164763846
$ws.Range("D7").NumberFormat = "0"
$ws.Range("E7").Value = "jlozano@suraelec.com"

# --- Hyperlinks for the new email cells (creates the mailto relationships) ---

$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:ravila@suraelec.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:ccharris@suraelec.com")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:jlozano@suraelec.com")

# Hyperlinks.Add forces its own style on the cell; put the normal Hyperlink
# cell style back on all 3 so they match the other email cells (E2:E4).
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"
$ws.Range("E7").Style = "Hyperlink"

# --- Expand Table1 to cover the new rows ---

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F7"))

# --- Column widths (best-fit for the new/affected columns) ---

$ws.Columns.Item(2).ColumnWidth = 14.83
$ws.Columns.Item(5).ColumnWidth = 21

# --- Selection cursor, matching the saved selection in the workbook ---

$ws.Range("D18").Select()
